$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: name + week number ---
$ws.Range("C1").Value = "Jesse Hare"
$ws.Range("E1").Value = 10

# --- Task rows (Stage / Task / Estimated / HoursSpent) ---
$ws.Range("A3").Value = "Project Build"
$ws.Range("B3").Value = "Work on iteration"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 3

$ws.Range("A4").Value = "Project Build"
$ws.Range("B4").Value = "Complete iteration"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 7

$ws.Range("A5").Value = "Proj analysis/elicitation"
$ws.Range("B5").Value = "Interview with client, redefine requirements from feedback"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 7

$ws.Range("A6").Value = "Project Build"
$ws.Range("B6").Value = "Work on final iteration"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 3

# --- Cumulative total label (D14 formula recalculates automatically) ---
$ws.Range("A14").Value = "Cumulative Total: 200"

# --- Column A got wider to fit the new text ---
$ws.Columns.Item(1).ColumnWidth = 23.5703125

# --- Selection moved ---
$ws.Range("D8").Select()
